# New scraping round: the "Mar 2014" total (D2) was re-scraped and came
# back slightly lower (12775 -> 12754). Update the figure, highlight the
# corrected cell in yellow so it's easy to spot, and leave the selection
# sitting on the cell that changed. All dependent formulas (F2, C3/E3,
# D9, E5, B13, B14, ...) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 12754

# Highlight the edited cell with a solid yellow fill (RGB 255,255,0).
$ws.Range("D2").Interior.Color = 65535

# Mirror the author's selection ending up on the corrected cell.
$ws.Range("D2").Select()
